# Automatische test-sync: 2025-08-03 14:20:50
#
# - Appends a new "Planning / Afspraak" test-mail row to the Logs sheet.
# - Extends the conditional-formatting ranges on Logs to cover the new row.
# - Appends the matching aggregate row to the Dashboard sheet.
# - Extends the bar chart's category/value series references to include it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 8 with the new test-mail entry.
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A8").Value = "Kun jij dit even regelen?"
$logs.Range("B8").Value = "mailmind.test@zohomail.eu"
$logs.Range("C8").Value = "Testmail #1: Kun jij dit even regelen?"
$logs.Range("D8").Value = "Planning / Afspraak"
$logs.Range("E8").Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$logs.Range("F8").Value = "2025-08-03 14:19:58"
$logs.Range("G8").Value = "Ja"
$logs.Range("H8").Value = "Ja"
$logs.Range("I8").Value = "Nee"
$logs.Range("J8").Value = "Nee"

# Extend the existing conditional formats (previously ...2:...7) to cover row 8.
$logs.Range("D2:D7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D8"))
$logs.Range("G2:G7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G8"))
$logs.Range("H2:H7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H8"))
$logs.Range("I2:I7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I8"))
$logs.Range("J2:J7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J8"))

# ---------------------------------------------------------------------
# 2. Dashboard sheet: append the aggregate row for the new category.
# ---------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A3").Value = "Planning / Afspraak"
$dashboard.Range("B3").Value = 1

# ---------------------------------------------------------------------
# 3. Chart: widen the category/value series to include the new row.
# ---------------------------------------------------------------------
$chartObj = $dashboard.ChartObjects().Item(1)
$series = $chartObj.Chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$3,Dashboard!`$B`$2:`$B`$3,1)"
